$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 98
$ws.Range("B3").Value = 82
$ws.Range("C3").Value = 82
$ws.Range("F3").Value = 146
$ws.Range("I3").Value = 199
$ws.Range("B6").Value = 395
$ws.Range("C6").Value = 508
$ws.Range("D6").Value = 443
$ws.Range("E6").Value = 511
$ws.Range("F6").Value = 579
$ws.Range("H6").Value = 480
$ws.Range("I6").Value = 518
$ws.Range("B7").Value = 535
$ws.Range("C7").Value = 670
$ws.Range("D7").Value = 686
$ws.Range("E7").Value = 749
$ws.Range("F7").Value = 833
$ws.Range("H7").Value = 769
$ws.Range("I7").Value = 863

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F3").Value = 5
$ws.Range("E6").Value = 56
$ws.Range("E7").Value = 69
$ws.Range("F7").Value = 64

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("F5").Value = 14
$ws.Range("F6").Value = 18

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("B4").Value = 13
$ws.Range("D4").Value = 5
$ws.Range("B5").Value = 17
$ws.Range("D5").Value = 6

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("B6").Value = 34
$ws.Range("B7").Value = 40

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 15

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F5").Value = 18
$ws.Range("C19").Value = 11
$ws.Range("F19").Value = 26
$ws.Range("E20").Value = 5
$ws.Range("B28").Value = 40
$ws.Range("E32").Value = 69
$ws.Range("F32").Value = 64
$ws.Range("B45").Value = 4
$ws.Range("C53").Value = 62
$ws.Range("D53").Value = 79
$ws.Range("F53").Value = 88
$ws.Range("I53").Value = 128
$ws.Range("H61").Value = 9
$ws.Range("I61").Value = 4
$ws.Range("E65").Value = 20
$ws.Range("F74").Value = 12
$ws.Range("D79").Value = 7
$ws.Range("B80").Value = 17
$ws.Range("D80").Value = 6
$ws.Range("B82").Value = 15
$ws.Range("B98").Value = 535
$ws.Range("C98").Value = 670
$ws.Range("D98").Value = 686
$ws.Range("E98").Value = 749
$ws.Range("F98").Value = 833
$ws.Range("H98").Value = 769
$ws.Range("I98").Value = 863

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F2").Value = 8
$ws.Range("I3").Value = 31
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 49
$ws.Range("C7").Value = 62
$ws.Range("D7").Value = 79
$ws.Range("F7").Value = 88
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 7

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("E5").Value = 17
$ws.Range("E6").Value = 20

$ws = $wb.Worksheets.Item('River North')
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 12

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("C3").Value = 1
$ws.Range("F6").Value = 19
$ws.Range("C7").Value = 11
$ws.Range("F7").Value = 26

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("B3").Value = 1
$ws.Range("B6").Value = 4

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
